# "Commiting documents for scrum flavor"
# Rebrand the "Goblins" team sheets/charts to "Team", and the
# "Originations" Quality-Gates bug-trend artifacts to "Product".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the three team sheets (Quality Gates & Sheet1 untouched)
# ---------------------------------------------------------------
$wsBurnup    = $wb.Worksheets.Item("Goblins Release Burn-up")
$wsVelocity  = $wb.Worksheets.Item("Goblins Velocity")
$wsFocus     = $wb.Worksheets.Item("Goblins Focus Factor")
$wsQuality   = $wb.Worksheets.Item("Quality Gates")

$wsBurnup.Name   = "Team Release Burn-up"
$wsVelocity.Name = "Team Velocity"
$wsFocus.Name    = "Team Focus Factor"

# ---------------------------------------------------------------
# 2. Quality Gates A2 label: Originations -> Product bugs trend
# ---------------------------------------------------------------
$wsQuality.Range("A2").Value = "2018R1 - Product Bugs Trend (end of each sprint)"

# ---------------------------------------------------------------
# 3. Release Burn-up chart: title + all series source formulas
# ---------------------------------------------------------------
$chartBurnup = $wsBurnup.ChartObjects().Item(1).Chart
$chartBurnup.ChartTitle.Text = "Team 2018R1 Release Burn Up"
for ($i = 1; $i -le $chartBurnup.SeriesCollection().Count; $i++) {
  $ser = $chartBurnup.SeriesCollection().Item($i)
  $ser.Formula = $ser.Formula.Replace("'Goblins Release Burn-up'", "'Team Release Burn-up'")
}

# ---------------------------------------------------------------
# 4. Velocity chart: title + all series source formulas
# ---------------------------------------------------------------
$chartVelocity = $wsVelocity.ChartObjects().Item(1).Chart
$chartVelocity.ChartTitle.Text = "Team 2018R1 Sprint Velocity & Average Velocity"
for ($i = 1; $i -le $chartVelocity.SeriesCollection().Count; $i++) {
  $ser = $chartVelocity.SeriesCollection().Item($i)
  $ser.Formula = $ser.Formula.Replace("'Goblins Velocity'", "'Team Velocity'")
}

# ---------------------------------------------------------------
# 5. Focus Factor chart: title + series source formula
# ---------------------------------------------------------------
$chartFocus = $wsFocus.ChartObjects().Item(1).Chart
$chartFocus.ChartTitle.Text = "Team Focus Factor"
for ($i = 1; $i -le $chartFocus.SeriesCollection().Count; $i++) {
  $ser = $chartFocus.SeriesCollection().Item($i)
  $ser.Formula = $ser.Formula.Replace("'Goblins Focus Factor'", "'Team Focus Factor'")
}

# ---------------------------------------------------------------
# 6. Quality Gates "Bugs Trend" chart title: Originations -> Product
#    (the sibling "Veracode Medium Flaws Trend" chart is unchanged)
# ---------------------------------------------------------------
$chartBugsTrend = $wsQuality.ChartObjects().Item(1).Chart
$chartBugsTrend.ChartTitle.Text = "Product Bugs Trend @end of each sprint"

# ---------------------------------------------------------------
# 7. Restore each sheet's selection, then finish with the
#    Release Burn-up sheet active/selected (matches saved view state)
# ---------------------------------------------------------------
$wsQuality.Range("E17").Select()
$wsFocus.Range("J22").Select()
$wsVelocity.Range("P34").Select()
$wsBurnup.Activate()
$wsBurnup.Range("H24").Select()
